$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Latitude / longitude for existing rows 2-5
$ws.Range("D2").Value = 48.86614866335567
$ws.Range("E2").Value = 2.31900834175702

$ws.Range("D3").Value = 48.8554814
$ws.Range("E3").Value = 2.3604077

$ws.Range("D4").Value = 45.7588923
$ws.Range("E4").Value = 4.8309221

$ws.Range("D5").Value = 45.7640318
$ws.Range("E5").Value = 4.8356904

# New row 6
$ws.Range("A6").Value = "HATIER"
$ws.Range("B6").Value = "Cléo"
$ws.Range("C6").Value = "Promenade des Anglais, 06000 Nice"
$ws.Range("D6").Value = 43.6859892
$ws.Range("E6").Value = 7.237476542487647
